$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3786.6667
$ws.Range("I74").Value = 3390.9092
$ws.Range("J74").Value = 4015.7896
$ws.Range("K74").Value = 3390.9092
$ws.Range("L74").Value = 4015.7896
$ws.Range("M74").Value = -2454.9092
$ws.Range("N74").Value = -5887.7896

$ws.Range("H76").Value = 5138.722
$ws.Range("I76").Value = 3642.4285
$ws.Range("J76").Value = 6090.909
$ws.Range("K76").Value = 3642.4285
$ws.Range("L76").Value = 6090.909
$ws.Range("M76").Value = -3327.4285
$ws.Range("N76").Value = -6720.909

$ws.Range("H77").Value = 3786.6667
$ws.Range("I77").Value = 3390.9092
$ws.Range("J77").Value = 4015.7896
$ws.Range("K77").Value = 16954.546
$ws.Range("L77").Value = 20078.948
$ws.Range("M77").Value = -12274.546
$ws.Range("N77").Value = -29438.948

$ws.Range("H79").Value = 5138.722
$ws.Range("I79").Value = 3642.4285
$ws.Range("J79").Value = 6090.909
$ws.Range("K79").Value = 3642.4285
$ws.Range("L79").Value = 6090.909
$ws.Range("M79").Value = -2550.4285
$ws.Range("N79").Value = -8274.909

$ws.Range("H86").Value = 3968.6667
$ws.Range("I86").Value = 3968.6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3968.6667
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2845.6667
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3968.6667
$ws.Range("I89").Value = 3968.6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 19843.3335
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -14227.3335
$ws.Range("N89").ClearContents()

$ws.Range("H116").Value = 3222.6667
$ws.Range("I116").Value = 2629.75
$ws.Range("J116").Value = 3519.125
$ws.Range("K116").Value = 2629.75
$ws.Range("L116").Value = 3519.125
$ws.Range("M116").Value = 812.25
$ws.Range("N116").Value = -10403.125


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1083.6111
$ws.Range("I2").Value = 1077.5714
$ws.Range("K2").Value = 1077.5714
$ws.Range("M2").Value = -964.5714

$ws.Range("H95").Value = 39331
$ws.Range("J95").Value = 39331
$ws.Range("L95").Value = 39331
$ws.Range("N95").Value = -44823

$ws.Range("H97").Value = 6907.2104
$ws.Range("I97").Value = 7474.0713
$ws.Range("J97").Value = 5320
$ws.Range("K97").Value = 7474.0713
$ws.Range("L97").Value = 5320
$ws.Range("M97").Value = -6978.0713
$ws.Range("N97").Value = -6312

$ws.Range("H116").Value = 1083.6111
$ws.Range("I116").Value = 1077.5714
$ws.Range("K116").Value = 1077.5714
$ws.Range("M116").Value = 1216.4286


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1083.6111
$ws.Range("I3").Value = 1077.5714
$ws.Range("K3").Value = 1077.5714
$ws.Range("M3").Value = -963.5714

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 9656
$ws.Range("J95").Value = 9656
$ws.Range("L95").Value = 9656
$ws.Range("N95").Value = -15148

$ws.Range("H105").Value = 1699.8334
$ws.Range("I105").Value = 1699.8334
$ws.Range("K105").Value = 1699.8334
$ws.Range("M105").Value = 47.16660000000002


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -64118

$ws.Range("H114").Value = 1518.4
$ws.Range("I114").Value = 1193.6364
$ws.Range("J114").Value = 1915.3334
$ws.Range("K114").Value = 3580.9092
$ws.Range("L114").Value = 5746.0002
$ws.Range("M114").Value = -326.9092000000001
$ws.Range("N114").Value = -12254.0002

$ws.Range("H132").Value = 2674
$ws.Range("I132").Value = 796.6667
$ws.Range("J132").Value = 5490
$ws.Range("K132").Value = 7170.0003
$ws.Range("L132").Value = 49410
$ws.Range("M132").Value = -4640.0003
$ws.Range("N132").Value = -54470


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2500300
$ws.Range("J3").Value = 600
$ws.Range("L3").Value = 600
$ws.Range("N3").Value = -832

$ws.Range("H4").Value = 1580.6
$ws.Range("J4").Value = 3050
$ws.Range("L4").Value = 3050
$ws.Range("N4").Value = -3274

$ws.Range("H80").Value = 13892265
$ws.Range("I80").Value = 27780030
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 27780030
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -27779032
$ws.Range("N80").Value = -6496

$ws.Range("H83").Value = 13892265
$ws.Range("I83").Value = 27780030
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 138900150
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -138895158
$ws.Range("N83").Value = -32484

$ws.Range("H132").Value = 3936.689
$ws.Range("I132").Value = 2707.5
$ws.Range("J132").Value = 6962.385
$ws.Range("K132").Value = 8122.5
$ws.Range("L132").Value = 20887.155
$ws.Range("M132").Value = -5592.5
$ws.Range("N132").Value = -25947.155


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4418.787
$ws.Range("J7").Value = 5334.1113
$ws.Range("L7").Value = 5334.1113
$ws.Range("N7").Value = -5558.1113

$ws.Range("H40").Value = 6033.6665
$ws.Range("I40").Value = 5471.857
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 5471.857
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -5335.857
$ws.Range("N40").Value = -8272

$ws.Range("H95").Value = 19343
$ws.Range("J95").Value = 19343
$ws.Range("L95").Value = 19343
$ws.Range("N95").Value = -24835

$ws.Range("H96").Value = 19194.428
$ws.Range("I96").Value = 13179
$ws.Range("K96").Value = 13179
$ws.Range("M96").Value = -10433

$ws.Range("H126").Value = 4418.787
$ws.Range("J126").Value = 5334.1113
$ws.Range("L126").Value = 16002.3339
$ws.Range("N126").Value = -20942.3339


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 6669333.5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 6669333.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 6669333.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -6669557.5

$ws.Range("H122").Value = 1688.2916
$ws.Range("I122").Value = 1816.8422
$ws.Range("K122").Value = 5450.5266
$ws.Range("M122").Value = -3000.5266

